# Weekly data refresh: insert a new observation row at the top of the data
# block (row 10) and push the existing rows down by one. Excel's native
# Rows.Insert() shifts rows 10:121 down to 11:122 (carrying all of their
# data/formatting with them), which is exactly the row-shift pattern seen
# in the target diff. We then just need to populate the freshly-inserted
# row 10 with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 10; rows 10-121 shift down to 11-122.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with this week's observation.
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = 45043
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 100112030
$ws.Cells.Item(10, 7).Value = "Poroto granado"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 30
$ws.Cells.Item(10, 11).Value = 35000
$ws.Cells.Item(10, 12).Value = 35000
$ws.Cells.Item(10, 13).Value = 35000
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 1400
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
